$d = $word.ActiveDocument
$d.Content.Find.Execute("juse ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "use ", 2)
